$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6289.0454
$ws.Range("I113").Value = 7446.7334
$ws.Range("K113").Value = 7446.7334
$ws.Range("M113").Value = -4192.7334

$ws.Range("H116").Value = 11048.583
$ws.Range("I116").Value = 12058.5
$ws.Range("K116").Value = 12058.5
$ws.Range("M116").Value = -8616.5

$ws.Range("H138").Value = 1865.3163
$ws.Range("J138").Value = 2507.3547
$ws.Range("L138").Value = 7522.0641
$ws.Range("N138").Value = -17802.0641

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 10056.75
$ws.Range("I14").Value = 7185.3335
$ws.Range("J14").Value = 11779.6
$ws.Range("K14").Value = 7185.3335
$ws.Range("L14").Value = 11779.6
$ws.Range("M14").Value = -7010.3335
$ws.Range("N14").Value = -12129.6

$ws.Range("H30").Value = 4691.6665
$ws.Range("I30").Value = 1245
$ws.Range("J30").Value = 11585
$ws.Range("K30").Value = 1245
$ws.Range("L30").Value = 11585
$ws.Range("M30").Value = -1095
$ws.Range("N30").Value = -11885

$ws.Range("H61").Value = 4576
$ws.Range("I61").Value = 2628.2144
$ws.Range("K61").Value = 2628.2144
$ws.Range("M61").Value = -2416.2144

$ws.Range("H74").Value = 3691.4167
$ws.Range("I74").Value = 3316.484
$ws.Range("J74").Value = 6016
$ws.Range("K74").Value = 3316.484
$ws.Range("L74").Value = 6016
$ws.Range("M74").Value = -2442.484
$ws.Range("N74").Value = -7764

$ws.Range("H77").Value = 3691.4167
$ws.Range("I77").Value = 3316.484
$ws.Range("J77").Value = 6016
$ws.Range("K77").Value = 16582.42
$ws.Range("L77").Value = 30080
$ws.Range("M77").Value = -12214.42
$ws.Range("N77").Value = -38816

$ws.Range("H97").Value = 1153.7273
$ws.Range("I97").Value = 925
$ws.Range("J97").Value = 2183
$ws.Range("K97").Value = 925
$ws.Range("L97").Value = 2183
$ws.Range("M97").Value = -429
$ws.Range("N97").Value = -3175

$ws.Range("H119").Value = 40416.5
$ws.Range("J119").Value = 40416.5
$ws.Range("L119").Value = 40416.5
$ws.Range("N119").Value = -50092.5

$ws.Range("H122").Value = 4596.125
$ws.Range("I122").Value = 3098
$ws.Range("K122").Value = 9294
$ws.Range("M122").Value = -6844

$ws.Range("H136").Value = 4576
$ws.Range("I136").Value = 2628.2144
$ws.Range("K136").Value = 7884.6432
$ws.Range("M136").Value = -5334.6432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 20016
$ws.Range("I29").Value = 20016
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 20016
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -19727
$ws.Range("N29").ClearContents()

$ws.Range("H94").Value = 518
$ws.Range("I94").Value = 532.75
$ws.Range("K94").Value = 532.75
$ws.Range("M94").Value = -81.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2687.875
$ws.Range("I31").Value = 2282
$ws.Range("K31").Value = 2282
$ws.Range("M31").Value = -1987

$ws.Range("H34").Value = 2687.875
$ws.Range("I34").Value = 2282
$ws.Range("K34").Value = 2282
$ws.Range("M34").Value = -2080

$ws.Range("H103").Value = 30000
$ws.Range("I103").Value = 30000
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 30000
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -28828
$ws.Range("N103").ClearContents()

$ws.Range("H122").Value = 1056.8572
$ws.Range("I122").Value = 781.3333
$ws.Range("K122").Value = 2343.9999
$ws.Range("M122").Value = 106.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 400.875
$ws.Range("I97").Value = 371.6
$ws.Range("J97").Value = 449.66666
$ws.Range("K97").Value = 1114.8
$ws.Range("L97").Value = 1348.99998
$ws.Range("M97").Value = -618.8000000000002
$ws.Range("N97").Value = -2340.99998

$ws.Range("H117").Value = 1403.375
$ws.Range("I117").Value = 357.8
$ws.Range("J117").Value = 1878.6364
$ws.Range("K117").Value = 1073.4
$ws.Range("L117").Value = 5635.9092
$ws.Range("M117").Value = 2368.6
$ws.Range("N117").Value = -12519.9092

$ws.Range("H121").Value = 90911440
$ws.Range("J121").Value = 5665.6665
$ws.Range("L121").Value = 16996.9995
$ws.Range("N121").Value = -19616.9995

$ws.Range("H122").Value = 929.1
$ws.Range("J122").Value = 1080.5714
$ws.Range("L122").Value = 9725.142600000001
$ws.Range("N122").Value = -14625.1426

$ws.Range("H136").Value = 3303.2932
$ws.Range("J136").Value = 5461.3076
$ws.Range("L136").Value = 16383.9228
$ws.Range("N136").Value = -26583.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 13781.5
$ws.Range("I3").Value = 8619.875
$ws.Range("K3").Value = 8619.875
$ws.Range("M3").Value = -8503.875

$ws.Range("H53").Value = 19000
$ws.Range("I53").Value = 19000
$ws.Range("K53").Value = 19000
$ws.Range("M53").Value = -18369

$ws.Range("H102").Value = 5810.1934
$ws.Range("I102").Value = 5781.4136
$ws.Range("K102").Value = 5781.4136
$ws.Range("M102").Value = -4159.4136

$ws.Range("H122").Value = 1427.4736
$ws.Range("I122").Value = 1358.6875
$ws.Range("K122").Value = 4076.0625
$ws.Range("M122").Value = -1626.0625

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7224.8335
$ws.Range("I40").Value = 7548.864
$ws.Range("K40").Value = 7548.864
$ws.Range("M40").Value = -7412.864

$ws.Range("H61").Value = 1011.6087
$ws.Range("I61").Value = 966.4545000000001
$ws.Range("K61").Value = 966.4545000000001
$ws.Range("M61").Value = -764.4545000000001

$ws.Range("H106").Value = 18115.75
$ws.Range("J106").Value = 18115.75
$ws.Range("L106").Value = 18115.75
$ws.Range("N106").Value = -20639.75

$ws.Range("H113").Value = 1011.6087
$ws.Range("I113").Value = 966.4545000000001
$ws.Range("K113").Value = 966.4545000000001
$ws.Range("M113").Value = 1203.5455

$ws.Range("H122").Value = 4854.7754
$ws.Range("I122").Value = 4514.857
$ws.Range("J122").Value = 6894.2856
$ws.Range("K122").Value = 13544.571
$ws.Range("L122").Value = 20682.8568
$ws.Range("M122").Value = -11094.571
$ws.Range("N122").Value = -25582.8568

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H136").Value = 2252.423
$ws.Range("I136").Value = 1687.6097
$ws.Range("K136").Value = 5062.8291
$ws.Range("M136").Value = -2512.8291

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 24639.4
$ws.Range("J41").Value = 24639.4
$ws.Range("L41").Value = 24639.4
$ws.Range("N41").Value = -25419.4

$ws.Range("H97").Value = 20999.75
$ws.Range("J97").Value = 20999.75
$ws.Range("L97").Value = 20999.75
$ws.Range("N97").Value = -22981.75

$ws.Range("H122").Value = 816
$ws.Range("I122").Value = 802
$ws.Range("J122").Value = 998
$ws.Range("K122").Value = 2406
$ws.Range("L122").Value = 2994
$ws.Range("M122").Value = 44
$ws.Range("N122").Value = -7894

$ws.Range("H132").Value = 2185.6182
$ws.Range("I132").Value = 1874.6364
$ws.Range("K132").Value = 5623.9092
$ws.Range("M132").Value = -3093.9092

$ws.Range("H136").Value = 2674.1343
$ws.Range("I136").Value = 2391.2068
$ws.Range("K136").Value = 7173.6204
$ws.Range("M136").Value = -4623.6204
